$d = $word.ActiveDocument

# 1. Remove the now-redundant "Type resource" heading paragraph
#    (its table is kept, it simply moves up under the renamed first heading).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Type resource") {
        $p.Range.Delete()
        break
    }
}

# 2. Remove the first table entirely (the old "resourceRequest" summary table).
$d.Tables.Item(1).Delete()

# 3. Rename the first heading to reflect the merged content.
$d.Content.Find.Execute("Objet geolocalisation", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Objet geoResourceRequest", 2)

Write-Output "done"
